$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accelerometer samples to insert right after the header row (row 1),
# pushing the existing data down.
$newData = @(
    @(-0.9284301400184631, 1.425136804580689, -0.2127189040184021),
    @(-0.9707106351852418, 1.45836865901947, -0.1764526814222335),
    @(-1.115207254886627, 1.46594226360321, -0.1345747746527196),
    @(-1.505423545837401, 1.455123424530029, -0.2340321838855745),
    @(-1.116380929946899, 1.48697829246521, -0.4328413642942907),
    @(-1.109515905380249, 1.432106614112854, -0.3912773653864859),
    @(-1.141456544399262, 1.384602665901184, -0.2541450988501308),
    @(-1.13429856300354, 1.39785385131836, -0.2251825407147409),
    @(-1.063723325729371, 1.418689608573914, -0.2562501281499862),
    @(-1.018438935279846, 1.406062006950379, -0.2269966453313826),
    @(-1.137969434261322, 1.409385621547699, -0.1802991181612014)
)

$insertCount = $newData.Count

# Insert enough blank rows right below the header to hold the new samples.
$ws.Rows("2:$($insertCount + 1)").Insert()

# Row insertion copies formatting from the row above (the bold header);
# strip that back off so the new data rows stay unstyled like the rest.
$ws.Range("A2:C$($insertCount + 1)").ClearFormats()

# Fill the newly inserted rows with the new sample values.
for ($i = 0; $i -lt $insertCount; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $newData[$i][0]
    $ws.Cells.Item($row, 2).Value = $newData[$i][1]
    $ws.Cells.Item($row, 3).Value = $newData[$i][2]
}

# The previous last data row (old row 21, the final sample) was dropped
# from the dataset; after the insert it now lives at row 32.
$ws.Rows("32:32").Delete()
